$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B99").Value = "no"
$ws.Range("D101").Value = "event_log4.txt(64-94)"
$ws.Range("D117").Value = "event_log3.txt(87-98), event_log4.txt(29-46, 108-119)"
$ws.Range("D118").Value = "event_log2(28-44), event_log3.txt(100-114)"

$ws.Range("B122").Value = "no"
$ws.Range("B123").Value = "no"
$ws.Range("B124").Value = "no"
$ws.Range("B125").Value = "no"
$ws.Range("B126").Value = "no"
$ws.Range("B127").Value = "no"

$ws.Range("C104").Select()
